$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new weekly record as row 145, pushing the existing rows
# (old 145..187) down to (146..188).
$ws.Rows.Item(145).Insert()

$ws.Range("A145").Value = 11
$ws.Range("B145").Value = "Vega Monumental Concepción"
$ws.Range("C145").Value = "Bíobío"
$ws.Range("D145").Value = 44754
$ws.Range("E145").Value = 8
$ws.Range("F145").Value = 100112003
$ws.Range("G145").Value = "Ajo"
$ws.Range("H145").Value = "Chino"
$ws.Range("I145").Value = "Primera"
$ws.Range("J145").Value = 400
$ws.Range("K145").Value = 21000
$ws.Range("L145").Value = 22000
$ws.Range("M145").Value = 21500
$ws.Range("N145").Value = "`$/caja 10 kilos"
$ws.Range("O145").Value = "China"
$ws.Range("P145").Value = 2150
$ws.Range("Q145").Value = 10
$ws.Range("R145").Value = "Hortaliza"
